$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "야드" (yard) column header
$ws.Range("D1").Value = "야드"

# Rows that get yard value 2 (visible after filtering)
$d2Ranges = @(
    @(2,3), @(27,31), @(35,40), @(43,43), @(46,46), @(50,50), @(59,59), @(62,62), @(66,67), @(83,84)
)
foreach ($range in $d2Ranges) {
    for ($r = $range[0]; $r -le $range[1]; $r++) {
        $ws.Cells.Item($r, 4).Value = 2
    }
}

# Rows that get yard value 1 (hidden after filtering)
$d1Ranges = @(
    @(4,26), @(32,34), @(41,42), @(44,45), @(47,49), @(51,58), @(60,61), @(63,65), @(68,82)
)
foreach ($range in $d1Ranges) {
    for ($r = $range[0]; $r -le $range[1]; $r++) {
        $ws.Cells.Item($r, 4).Value = 1
    }
}

# Apply an AutoFilter over A1:D84 restricted to yard == 2
$ws.Range("A1:D84").AutoFilter(4, @("2"), 7)

# AutoFilter leaves behind the internal _FilterDatabase defined name
$filterDbName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet2!`$A`$1:`$D`$84")
$filterDbName.Visible = $false

# Match the cell selection left behind by the editing session
$ws.Range("G35").Select()
